$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "[Euclides-Tec. Soldagem-3B, Mayra-T. NãoMetalicos-3B, Euclides-Tec. Soldagem-3B, Mayra-T. NãoMetalicos-3B]"
$ws.Range("C3").Value = "[Rogério-Retífica-3B, Rogério-Retífica-3B, Rogério-Retífica-3B, Rogério-Retífica-3B]"
$ws.Range("D3").Value = "[Victor S.-Usin. CNC-3B, Victor S.-Usin. CNC-3B, Victor S.-Usin. CNC-3B, Victor S.-Usin. CNC-3B]"
$ws.Range("F3").Value = "Nilton-Elem. Máquinas"

$ws.Range("B4").Value = "[Euclides-Tec. Soldagem-3B, Gisele-Ens. Dest. Não Desti.-3B, Euclides-Tec. Soldagem-3B, Gisele-Ens. Dest. Não Desti.-3B]"
$ws.Range("C4").Value = "[Joel L.-Tec. Fundição-3B, Leandro-Mec. Manut. Equip. Ind-3B, Leandro-Mec. Manut. Equip. Ind-3B, Joel L.-Tec. Fundição-3B]"
$ws.Range("D4").Value = "[Aderci-Fresagem-3B, Aderci-Fresagem-3B, Aselmo-M. Motor Endot.-3B, Aselmo-M. Motor Endot.-3B]"
$ws.Range("F4").Value = "Nilton-Elem. Máquinas"

$ws.Range("B6").Value = "[Humberto-Coman. Pneumáticos-3B, Gisele-Ens. Dest. Não Desti.-3B, Humberto-Coman. Pneumáticos-3B, Gisele-Ens. Dest. Não Desti.-3B]"
$ws.Range("C6").Value = "[Leandro-M. S. R. AR Cond.-3B, Leandro-Mec. Manut. Equip. Ind-3B, Leandro-Mec. Manut. Equip. Ind-3B, Leandro-M. S. R. AR Cond.-3B]"
$ws.Range("D6").Value = "[Aderci-Fresagem-3B, Aderci-Fresagem-3B, Valmir-Calderaria-3B, Valmir-Calderaria-3B]"
$ws.Range("E6").Value = "[Elcio D.-Cont. Lóg. Prog. CLP-3B, Elcio D.-Cont. Lóg. Prog. CLP-3B, Elcio D.-Cont. Lóg. Prog. CLP-3B, Elcio D.-Cont. Lóg. Prog. CLP-3B]"
$ws.Range("F6").Value = "[Aselmo-M. Motor Endot.-3B, Ivan-Trat. Termicos-3B, Ivan-Trat. Termicos-3B, Aselmo-M. Motor Endot.-3B]"

$ws.Range("B7").Value = "[Ludoff-Coman. Hidraulicos-3B, Ludoff-Coman. Hidraulicos-3B, Ludoff-Coman. Hidraulicos-3B, Ludoff-Coman. Hidraulicos-3B]"
$ws.Range("C7").Value = "[Leandro-M. S. R. AR Cond.-3B, Paulo Rob.-M. A. Comp. CAD / CAM-3B, Paulo Rob.-M. A. Comp. CAD / CAM-3B, Leandro-M. S. R. AR Cond.-3B]"
$ws.Range("D7").Value = "[Ismail-Metrologia 2-3B, Ismail-Metrologia 2-3B, Ismail-Metrologia 2-3B, Ismail-Metrologia 2-3B]"
$ws.Range("E7").Value = "[Ludoff-Coman. Pneumáticos-3B, Ludoff-Coman. Pneumáticos-3B, Ludoff-Coman. Pneumáticos-3B, Ludoff-Coman. Pneumáticos-3B]"
$ws.Range("F7").Value = "[-, Ivan-Trat. Termicos-3B, Ivan-Trat. Termicos-3B, -]"

$ws.Range("C8").Value = "[Valmir-Calderaria-3B, Paulo Rob.-M. A. Comp. CAD / CAM-3B, Paulo Rob.-M. A. Comp. CAD / CAM-3B, Valmir-Calderaria-3B]"
$ws.Range("E8").Value = "[Joel L.-Tec. Fundição-3B, Joel L.-Tec. Fundição-3B, Humberto-Coman. Pneumáticos-3B, Humberto-Coman. Pneumáticos-3B]"
